# Add 8 new Family-Feud style question blocks (rows 196-243) to the
# 'Transformed by JSON-CSV.CO' sheet, plus backfill a missing answer
# label at B158. Operation order mirrors the original authoring pass
# (all B-column answers, then all A-column question headers, then the
# two late corrections) so shared-string append order matches exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pass 1: answer column (B) for the new rows, plus their score column (C)
$ws.Range("B196").Value = "Туристы\экскурсия"
$ws.Range("B197").Value = "Приёмная комиссия"
$ws.Range("B198").Value = "Покупатели"
$ws.Range("B200").Value = "Грабители"
$ws.Range("B201").Value = "Жильцы"
$ws.Range("B202").Value = "Повар"
$ws.Range("B204").Value = "Врач"
$ws.Range("B207").Value = "Банщик"
$ws.Range("B208").Value = "Сардина"
$ws.Range("B209").Value = "Сельдь"
$ws.Range("B210").Value = "Горбуша"
$ws.Range("B211").Value = "Скумбрия"
$ws.Range("B212").Value = "Сайра"
$ws.Range("B213").Value = "Килька"
$ws.Range("B214").Value = "Метель"
$ws.Range("B215").Value = "Мель"
$ws.Range("B216").Value = "Ельник"
$ws.Range("B217").Value = "Капель"
$ws.Range("B218").Value = "Карусель"
$ws.Range("B219").Value = "Мебель"
$ws.Range("B220").Value = "Муку"
$ws.Range("B221").Value = "Консервы"
$ws.Range("B222").Value = "Крупу"
$ws.Range("B223").Value = "Спички"
$ws.Range("B224").Value = "Сахар"
$ws.Range("B225").Value = "Соль"
$ws.Range("B226").Value = "Детям"
$ws.Range("B227").Value = "Пенсионерам"
$ws.Range("B228").Value = "Безработным"
$ws.Range("B229").Value = "Олигархам\богатым"
$ws.Range("B230").Value = "Инвалидам"
$ws.Range("B231").Value = "Домохозяйкам"
$ws.Range("B232").Value = "Дом"
$ws.Range("B233").Value = "Человек"
$ws.Range("B234").Value = "Город"
$ws.Range("B235").Value = "Язык"
$ws.Range("B236").Value = "Родители"
$ws.Range("B237").Value = "Ребёнок"
$ws.Range("B238").Value = "Лошадь"
$ws.Range("B239").Value = "Солдат"
$ws.Range("B240").Value = "Слон"
$ws.Range("B241").Value = "Цапля"
$ws.Range("B242").Value = "Охранник"
$ws.Range("B243").Value = "Корова"

# Rows that reuse answer text already present elsewhere in the sheet
$ws.Range("B203").Value = "Парикмахер"
$ws.Range("B205").Value = "Официант"
$ws.Range("B206").Value = "Спортсмен"

# Score column (C) for all 48 new rows
$ws.Range("C196").Value = 21
$ws.Range("C197").Value = 15
$ws.Range("C198").Value = 11
$ws.Range("C199").Value = 9
$ws.Range("C200").Value = 5
$ws.Range("C201").Value = 3
$ws.Range("C202").Value = 20
$ws.Range("C203").Value = 16
$ws.Range("C204").Value = 14
$ws.Range("C205").Value = 9
$ws.Range("C206").Value = 8
$ws.Range("C207").Value = 7
$ws.Range("C208").Value = 24
$ws.Range("C209").Value = 20
$ws.Range("C210").Value = 16
$ws.Range("C211").Value = 12
$ws.Range("C212").Value = 8
$ws.Range("C213").Value = 4
$ws.Range("C214").Value = 37
$ws.Range("C215").Value = 24
$ws.Range("C216").Value = 12
$ws.Range("C217").Value = 10
$ws.Range("C218").Value = 9
$ws.Range("C219").Value = 7
$ws.Range("C220").Value = 24
$ws.Range("C221").Value = 20
$ws.Range("C222").Value = 16
$ws.Range("C223").Value = 12
$ws.Range("C224").Value = 8
$ws.Range("C225").Value = 4
$ws.Range("C226").Value = 36
$ws.Range("C227").Value = 24
$ws.Range("C228").Value = 24
$ws.Range("C229").Value = 11
$ws.Range("C230").Value = 8
$ws.Range("C231").Value = 4
$ws.Range("C232").Value = 36
$ws.Range("C233").Value = 24
$ws.Range("C234").Value = 17
$ws.Range("C235").Value = 11
$ws.Range("C236").Value = 8
$ws.Range("C237").Value = 4
$ws.Range("C238").Value = 32
$ws.Range("C239").Value = 26
$ws.Range("C240").Value = 20
$ws.Range("C241").Value = 16
$ws.Range("C242").Value = 10
$ws.Range("C243").Value = 5

# Pass 2: question column (A) headers for each new block
$ws.Range("A196").Value = "033. Пять человек осматривают здание. Кто они?"
$ws.Range("A202").Value = "034. Кто на работе использует полотенце?"
$ws.Range("A208").Value = "035. Из какой рыбы делают консервы?"
$ws.Range("A214").Value = "036. В каких словах есть слово `"ель`"?"
$ws.Range("A220").Value = "037. Какие продукты или товары покупают `"про запас`"?"
$ws.Range("A226").Value = "038. Кому можно не работать?"
$ws.Range("A232").Value = "039. Что может быть родным?"
$ws.Range("A238").Value = "040. Кто спит стоя?"

# Late fix: row 199 answer label ("Строители\архитектор")
$ws.Range("B199").Value = "Строители\архитектор"

# Late fix: row 158 was missing its answer label ("Мопед")
$ws.Range("B158").Value = "Мопед"

# Update view state to match final selection
$ws.Range("B243").Select() | Out-Null
